$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column F (dSF)
$values = @{
    3  = 4
    4  = -7
    5  = -1
    6  = -9
    7  = 1
    8  = -1
    9  = 2
    10 = -1
    11 = -4
    12 = -2
    14 = -1
    15 = -1
    16 = 6
    17 = -4
    18 = -3
    19 = 3
    20 = -2
    21 = -2
    22 = -2
    23 = 3
    24 = 4
    25 = 1
    26 = 0
    27 = -2
    28 = -1
    30 = -1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
